$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AI" (sheet1): append a new data row (row 10) for "Charles Wilkes"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("AI")

$ws1.Range("A10").Value = "Charles Wilkes"
$ws1.Range("B10").Value = 346
$ws1.Range("C10").Value = 5
$ws1.Range("D10").Value = 88
$ws1.Range("E10").Value = 7
$ws1.Range("F10").Value = 397
$ws1.Range("G10").Value = 66
$ws1.Range("H10").Value = 609
$ws1.Range("I10").Value = 51
$ws1.Range("J10").Value = 21
$ws1.Range("K10").Value = 28
$ws1.Range("L10").Value = 0
$ws1.Range("M10").Value = 0
$ws1.Range("N10").Value = 0
$ws1.Range("O10").Value = 0
$ws1.Range("P10").Value = 0
$ws1.Range("Q10").Value = 0
$ws1.Range("R10").Value = 0
$ws1.Range("S10").Value = 1500
$ws1.Range("T10").Value = 43
$ws1.Range("U10").Value = 13
$ws1.Range("V10").Value = 10
$ws1.Range("W10").Value = 12
$ws1.Range("X10").Value = 22
$ws1.Range("Y10").Value = 500
$ws1.Range("Z10").Value = 7
$ws1.Range("AA10").Value = 46
$ws1.Range("AB10").Value = 28
$ws1.Range("AC10").Value = 15
$ws1.Range("AD10").Value = 4
$ws1.Range("AE10").Value = 338
$ws1.Range("AF10").Value = 87
$ws1.Range("AG10").Value = 535
$ws1.Range("AH10").Value = 46
$ws1.Range("AI10").Value = 11
$ws1.Range("AJ10").Value = 73
$ws1.Range("AK10").Value = 355
$ws1.Range("AL10").Value = 85
$ws1.Range("AM10").Value = 1000
$ws1.Range("AN10").Value = 1
$ws1.Range("AO10").Value = 60
$ws1.Range("AP10").Value = 39
$ws1.Range("AQ10").Value = 971
$ws1.Range("AR10").Value = 7
$ws1.Range("AS10").Value = 10
$ws1.Range("AT10").Value = 64
$ws1.Range("AU10").Value = 18
$ws1.Range("AV10").Value = 18
$ws1.Range("AW10").Value = 0
$ws1.Range("AX10").Value = 0
$ws1.Range("AY10").Value = 0
$ws1.Range("AZ10").Value = 0
$ws1.Range("BA10").Value = 0
$ws1.Range("BB10").Value = 179
$ws1.Range("BC10").Value = 20
$ws1.Range("BD10").Value = 1800
$ws1.Range("BE10").Value = 10
$ws1.Range("BF10").Value = 10
$ws1.Range("BG10").Value = 10
$ws1.Range("BH10").Value = 10
$ws1.Range("BI10").Value = 60
$ws1.Range("BJ10").Value = 1800
$ws1.Range("BK10").Value = 10
$ws1.Range("BL10").Value = 10
$ws1.Range("BM10").Value = 10
$ws1.Range("BN10").Value = 10
$ws1.Range("BO10").Value = 60
$ws1.Range("BP10").Value = 3000
$ws1.Range("BQ10").Value = 29
$ws1.Range("BR10").Value = 3000
$ws1.Range("BS10").Value = 61
$ws1.Range("BT10").Value = 0
$ws1.Range("BU10").Value = 0
$ws1.Range("BV10").Value = 0
$ws1.Range("BW10").Value = 0

# ---------------------------------------------------------------------------
# Sheet "List2" (sheet2): fill in a few extra counters (B/C/E columns)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("List2")

$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 1

$ws2.Range("B3").Value = 1
$ws2.Range("E3").Value = 37

$ws2.Range("B4").Value = 1
$ws2.Range("E4").Value = 28

$ws2.Range("C5").Value = 2

# ---------------------------------------------------------------------------
# View state: active sheet moves from "List2" to "AI", selection on "AI"
# moves to A10 (the newly added row), and "List2" selection moves to C6.
# ---------------------------------------------------------------------------
$null = $ws2.Range("C6").Select()

$null = $ws1.Activate()
$null = $ws1.Range("A10").Select()
